$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell A45 with corrected timestamp
$ws.Range("A45").Value = 44358.76795423611

# Append new row 46 with retrieved data
$ws.Range("A46").Value = 44359.76822617463
$ws.Range("B46").Value = 77175
$ws.Range("C46").Value = 64899
$ws.Range("D46").Value = 3531
$ws.Range("E46").Value = 2105
$ws.Range("F46").Value = 1470
$ws.Range("G46").Value = 20375
$ws.Range("H46").Value = 1483
$ws.Range("I46").Value = 879
$ws.Range("J46").Value = 188

# Match styling of column A (date format) for new row
$ws.Range("A46").NumberFormat = $ws.Range("A45").NumberFormat
